$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Material Ambient Color" (row 34) and "Material Specular Intensity" (row 35)
# fields were removed from the material subset description. Deleting these two
# rows shifts everything below them up by two rows, which matches the rest of
# the documentation table automatically (shared strings get pruned on save).
$ws.Rows("34:35").Delete()

# The conditional formatting range tracked the old last row of the table
# (row 123); after removing two rows it should track the new last row (121).
$conditions = $ws.Range("C5:C123").FormatConditions
for ($i = 1; $i -le $conditions.Count; $i++) {
    $conditions.Item($i).ModifyAppliesToRange($ws.Range("C5:C121"))
}

# Reflect where the author's selection ended up after the edit.
$ws.Rows(34).Select()
